$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The schedule table was missing the Tuesday session (it had been dropped by the
# print-range bug) and the Friday session. Insert a new row for Tuesday right
# above the existing Wednesday row (pushing Wednesday down to row 7), then
# append Friday as a new row 8.

$ws.Rows("6").Insert()

$ws.Range("B6").Value = "Вторник"
$ws.Range("C6").Value = "11:45:00"
$ws.Range("D6").Value = "Базы данных - практика"
$ws.Range("E6").Value = "В100"

$ws.Range("B8").Value = "Пятница"
$ws.Range("C8").Value = "13:45:00"
$ws.Range("D8").Value = "Базы данных - практика"
$ws.Range("E8").Value = "В101"
